$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- sheet1 ---
$ws1.Range("A2:H11").ClearContents()
$ws1.Range("A2").Value = 'MCQ'
$ws1.Range("B2").Value = 'mcq Clone'
$ws1.Range("C2").Value = 'option shuffle mcq'
Set-TextValue $ws1.Range("D2") '10'
$ws1.Range("F2").Value = 'a b d e '
$ws1.Range("A3").Value = 'Subjective'
$ws1.Range("B3").Value = 'subjective Clone Clone'
$ws1.Range("C3").Value = 'asdf'
Set-TextValue $ws1.Range("D3") '10'
$ws1.Range("A4").Value = 'Coding'
$ws1.Range("B4").Value = 'abc test 2'
$ws1.Range("C4").Value = 'test test'
Set-TextValue $ws1.Range("D4") '20'
$ws1.Range("E4").Value = 'C C++ Java JavaScript Python 3 '
$ws1.Range("A5").Value = 'MQ'
$ws1.Range("B5").Value = 'mq script'
$ws1.Range("C5").Value = '<!DOCTYPE html>'
Set-TextValue $ws1.Range("D5") '12'
$ws1.Range("A6").Value = 'Subjective'
$ws1.Range("B6").Value = 'subjective'
$ws1.Range("C6").Value = 'asdf'
Set-TextValue $ws1.Range("D6") '10'
$ws1.Range("A7").Value = 'Subjective'
$ws1.Range("B7").Value = 'SUB'
$ws1.Range("C7").Value = 'SDFA'
Set-TextValue $ws1.Range("D7") '23'
$ws1.Range("A8").Value = 'Subjective'
$ws1.Range("B8").Value = 'keyword creation subjective'
$ws1.Range("C8").Value = 'fddvc'
Set-TextValue $ws1.Range("D8") '22'
$ws1.Range("A9").Value = 'Web'
$ws1.Range("B9").Value = '8 july web'
$ws1.Range("C9").Value = 'sdffghjghj'
Set-TextValue $ws1.Range("D9") '20'
$ws1.Range("A10").Value = 'Web'
$ws1.Range("B10").Value = 'keyword creation web'
$ws1.Range("C10").Value = 'werdft'
Set-TextValue $ws1.Range("D10") '1234'
$ws1.Range("A11").Value = 'Web'
$ws1.Range("B11").Value = 'daf'
$ws1.Range("C11").Value = 'ghj'
Set-TextValue $ws1.Range("D11") '23'

# --- sheet2 ---
$ws2.Range("A2:H12").ClearContents()
$ws2.Range("A2").Value = 'Subjective'
$ws2.Range("B2").Value = 'subjective'
$ws2.Range("C2").Value = 'asdf'
Set-TextValue $ws2.Range("D2") '10'
$ws2.Range("A3").Value = 'Subjective'
$ws2.Range("B3").Value = 'subjective file uplaod'
$ws2.Range("C3").Value = 'fghjkl'
Set-TextValue $ws2.Range("D3") '10'
$ws2.Range("A4").Value = 'Subjective'
$ws2.Range("B4").Value = 'subjective script'
$ws2.Range("C4").Value = '<!DOCTYPE html>'
Set-TextValue $ws2.Range("D4") '10'
$ws2.Range("A5").Value = 'Subjective'
$ws2.Range("B5").Value = 'SUB'
$ws2.Range("C5").Value = 'SDFA'
Set-TextValue $ws2.Range("D5") '23'
$ws2.Range("A6").Value = 'MCQ'
$ws2.Range("B6").Value = 'mnb key question'
$ws2.Range("C6").Value = 'hajkds'
Set-TextValue $ws2.Range("D6") '10'
$ws2.Range("F6").Value = 'a b '
$ws2.Range("A7").Value = 'Coding'
$ws2.Range("B7").Value = 'abc test 2'
$ws2.Range("C7").Value = 'test test'
Set-TextValue $ws2.Range("D7") '20'
$ws2.Range("E7").Value = 'C C++ Java JavaScript Python 3 '
$ws2.Range("A8").Value = 'MQ'
$ws2.Range("B8").Value = 'mq script'
$ws2.Range("C8").Value = '<!DOCTYPE html>'
Set-TextValue $ws2.Range("D8") '12'
$ws2.Range("A9").Value = 'MCQ'
$ws2.Range("B9").Value = 'mcqabc@'
$ws2.Range("C9").Value = 'ghajkld;f'
Set-TextValue $ws2.Range("D9") '10'
$ws2.Range("F9").Value = '<script>alert("Hello! This is a JavaScript alert.");</script> dfghj '
$ws2.Range("A10").Value = 'Web'
$ws2.Range("B10").Value = '8 july web'
$ws2.Range("C10").Value = 'sdffghjghj'
Set-TextValue $ws2.Range("D10") '20'
$ws2.Range("A11").Value = 'Web'
$ws2.Range("B11").Value = 'keyword creation web'
$ws2.Range("C11").Value = 'werdft'
Set-TextValue $ws2.Range("D11") '1234'
$ws2.Range("A12").Value = 'Web'
$ws2.Range("B12").Value = 'daf'
$ws2.Range("C12").Value = 'ghj'
Set-TextValue $ws2.Range("D12") '23'

# --- sheet3 ---
$ws3.Range("A2:H4").ClearContents()
$ws3.Range("A2").Value = 'Coding'
$ws3.Range("B2").Value = 'coding 1'
$ws3.Range("C2").Value = 'jksjfhksajdfhkas'
Set-TextValue $ws3.Range("D2") '20'
$ws3.Range("E2").Value = 'Java '
$ws3.Range("A3").Value = 'Coding'
$ws3.Range("B3").Value = 'coding 1 Clone'
$ws3.Range("C3").Value = 'jksjfhksajdfhkas'
Set-TextValue $ws3.Range("D3") '20'
$ws3.Range("E3").Value = 'Java '
$ws3.Range("A4").Value = 'Coding'
$ws3.Range("B4").Value = 'abc test 2'
$ws3.Range("C4").Value = 'test test'
Set-TextValue $ws3.Range("D4") '20'
$ws3.Range("E4").Value = 'C C++ Java JavaScript Python 3 '

